$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 2100
$ws.Range("I49").Value = 2100
$ws.Range("K49").Value = 6300
$ws.Range("M49").Value = -6164
$ws.Range("H53").Value = 35.42857
$ws.Range("I53").Value = 31.6
$ws.Range("K53").Value = 31.6
$ws.Range("M53").Value = 605.4
$ws.Range("H76").Value = 3000
$ws.Range("J76").Value = 3000
$ws.Range("L76").Value = 3000
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3000
$ws.Range("J79").Value = 3000
$ws.Range("L79").Value = 3000
$ws.Range("N79").Value = -5184
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H137").Value = 2139.8
$ws.Range("I137").Value = 2171.2144
$ws.Range("K137").Value = 6513.6432
$ws.Range("M137").Value = -3963.6432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2154.65
$ws.Range("I61").Value = 1439.5333
$ws.Range("J61").Value = 4300
$ws.Range("K61").Value = 1439.5333
$ws.Range("L61").Value = 4300
$ws.Range("M61").Value = -1227.5333
$ws.Range("N61").Value = -4724
$ws.Range("H76").Value = 33288
$ws.Range("J76").Value = 33288
$ws.Range("L76").Value = 33288
$ws.Range("N76").Value = -33964
$ws.Range("H79").Value = 33288
$ws.Range("J79").Value = 33288
$ws.Range("L79").Value = 33288
$ws.Range("N79").Value = -35628
$ws.Range("H122").Value = 4017.4167
$ws.Range("I122").Value = 4232.6665
$ws.Range("J122").Value = 3802.1667
$ws.Range("K122").Value = 12697.9995
$ws.Range("L122").Value = 11406.5001
$ws.Range("M122").Value = -10247.9995
$ws.Range("N122").Value = -16306.5001
$ws.Range("H136").Value = 2154.65
$ws.Range("I136").Value = 1439.5333
$ws.Range("J136").Value = 4300
$ws.Range("K136").Value = 4318.5999
$ws.Range("L136").Value = 12900
$ws.Range("M136").Value = -1768.5999
$ws.Range("N136").Value = -18000

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2943.5
$ws.Range("I86").Value = 1845.4
$ws.Range("J86").Value = 4773.6665
$ws.Range("K86").Value = 1845.4
$ws.Range("L86").Value = 4773.6665
$ws.Range("M86").Value = -722.4000000000001
$ws.Range("N86").Value = -7019.6665
$ws.Range("H89").Value = 2943.5
$ws.Range("I89").Value = 1845.4
$ws.Range("J89").Value = 4773.6665
$ws.Range("K89").Value = 9227
$ws.Range("L89").Value = 23868.3325
$ws.Range("M89").Value = -3611
$ws.Range("N89").Value = -35100.3325
$ws.Range("H105").Value = 2972
$ws.Range("I105").Value = 2963
$ws.Range("J105").Value = 2999
$ws.Range("K105").Value = 2963
$ws.Range("L105").Value = 2999
$ws.Range("M105").Value = -1216
$ws.Range("N105").Value = -6493
$ws.Range("H107").Value = 4998.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4998.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4998.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -8838.5
$ws.Range("H134").Value = 3055.5557
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3039.4
$ws.Range("I80").Value = 2932.6667
$ws.Range("J80").Value = 3199.5
$ws.Range("K80").Value = 2932.6667
$ws.Range("L80").Value = 3199.5
$ws.Range("M80").Value = -1934.6667
$ws.Range("N80").Value = -5195.5
$ws.Range("H83").Value = 3039.4
$ws.Range("I83").Value = 2932.6667
$ws.Range("J83").Value = 3199.5
$ws.Range("K83").Value = 14663.3335
$ws.Range("L83").Value = 15997.5
$ws.Range("M83").Value = -9671.333500000001
$ws.Range("N83").Value = -25981.5
$ws.Range("H97").Value = 354
$ws.Range("I97").Value = 343
$ws.Range("K97").Value = 343
$ws.Range("M97").Value = 153
$ws.Range("H107").Value = 1966.3846
$ws.Range("I107").Value = 1284
$ws.Range("J107").Value = 3058.2
$ws.Range("K107").Value = 1284
$ws.Range("L107").Value = 3058.2
$ws.Range("M107").Value = 636
$ws.Range("N107").Value = -6898.2
$ws.Range("H122").Value = 3023
$ws.Range("J122").Value = 1943
$ws.Range("L122").Value = 5829
$ws.Range("N122").Value = -10729
$ws.Range("H128").Value = 30768
$ws.Range("J128").Value = 30768
$ws.Range("L128").Value = 30768
$ws.Range("N128").Value = -40728

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1999.6666
$ws.Range("I22").Value = 1999
$ws.Range("K22").Value = 1999
$ws.Range("M22").Value = -1704
$ws.Range("H27").Value = 1999.6666
$ws.Range("I27").Value = 1999
$ws.Range("K27").Value = 1999
$ws.Range("M27").Value = -1892
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H62").Value = 1400
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376
$ws.Range("H65").Value = 1400
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880
$ws.Range("H81").Value = 3557.7693
$ws.Range("I81").Value = 3925.1
$ws.Range("J81").Value = 2333.3333
$ws.Range("K81").Value = 7850.2
$ws.Range("L81").Value = 4666.6666
$ws.Range("M81").Value = -6789.2
$ws.Range("N81").Value = -6788.6666
$ws.Range("H84").Value = 3557.7693
$ws.Range("I84").Value = 3925.1
$ws.Range("J84").Value = 2333.3333
$ws.Range("K84").Value = 39251
$ws.Range("L84").Value = 23333.333
$ws.Range("M84").Value = -33947
$ws.Range("N84").Value = -33941.333
$ws.Range("H107").Value = 2219.7778
$ws.Range("I107").Value = 1082.8334
$ws.Range("J107").Value = 4493.6665
$ws.Range("K107").Value = 3248.5002
$ws.Range("L107").Value = 13480.9995
$ws.Range("M107").Value = -1328.5002
$ws.Range("N107").Value = -17320.9995
$ws.Range("H132").Value = 1888.5143
$ws.Range("J132").Value = 4643.6665
$ws.Range("L132").Value = 13930.9995
$ws.Range("N132").Value = -18990.9995
$ws.Range("H136").Value = 1011.1739
$ws.Range("I136").Value = 1002.7143
$ws.Range("K136").Value = 3008.1429
$ws.Range("M136").Value = -458.1428999999998
